$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Keep only the 2008-2015 columns when transposing; 1998/2003 are dropped
# per the commit ("Rows transposed to columns").
$years  = @(2008, 2009, 2010, 2011, 2012, 2013, 2014, 2015)
$values = @(17873.635440000002, 15932.3683, 15554.491840000001, 15489.816060000001, 16890.79897, 18442.00488, 17835.548739999998, 18820.167510000003)

# Grab the existing "travels_private" label before the source range is
# wiped out, so the new header cell reuses the exact same text (incl. its
# leading non-breaking spaces) instead of retyping it.
$label = $ws.Range("A2").Value2

# Wipe the old row-oriented range (A1:K2), then drop the row-level
# height/format overrides that used to span columns A:K so the new rows
# don't inherit stale formatting from row 2.
$ws.Range("A1:K2").Clear()
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(2).EntireRow.AutoFit()

# Headers
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = $label

# Transposed data: years down column A, values down column B
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Re-apply the sheet's data style (right-aligned, vertically centered) to
# the new A1:B9 block.
$dataRange = $ws.Range("A1:B9")
$dataRange.Style = "Normal"
$dataRange.HorizontalAlignment = -4152  # xlRight
$dataRange.VerticalAlignment = -4108    # xlCenter

# Column B needs a width like column A already has; size it to its content.
$ws.Range("B1:B9").EntireColumn.AutoFit()

$ws.Rows("2:3").Select()
